# DailyWorkReport.xlsx edit
# - Update the existing D47 hours total from 4.5 -> 7
# - Append a new weekly block (rows 49-51) for 2025-01-20 (serial 45677),
#   mirroring the existing "Domm / Meeting / Study" row layout, followed
#   by a blank formatted row 52 (row 48 stays an empty gap row, matching
#   the pattern already used elsewhere in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Friday (2025-01-17) "Advance C# final demo" hours.
$ws.Range("D47").Value = 7

# Clone the formatting (borders / number formats / alignment) of the most
# recent complete day-block (rows 44:47) down onto the new block (49:52) so
# the new rows pick up the same styles used throughout the sheet.
$ws.Range("A44:D47").Copy()
$ws.Range("A49:D52").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 49: new day header (date + "Domm").
$ws.Range("A49").Value = 45677
$ws.Range("B49").Value = "Domm"
$ws.Range("D49").Value = 0.25

# Row 50: Meeting / General Discussion.
$ws.Range("B50").Value = "Meeting"
$ws.Range("C50").Value = "General Discussion"
$ws.Range("D50").Value = 0.25

# Row 51: Study / Advance C# final demo.
$ws.Range("B51").Value = "Study"
$ws.Range("C51").Value = "Advance C# final demo"
$ws.Range("D51").Value = 7.5

# Row 52 is left blank (just inherits the copied formatting) as the new
# trailing row, matching the diff.

# Move the active selection / scrolled viewport to the new bottom cell.
$ws.Range("D52").Select()
$excel.ActiveWindow.ScrollRow = 33
